$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 49, shifting existing rows 49-79 down to 50-80.
$ws.Rows.Item(49).Insert()

# Populate the newly inserted row 49 with the new weekly entry.
$ws.Range('A49').Value = 4
$ws.Range('B49').Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range('C49').Value = 'Los Lagos'
$ws.Range('D49').Value = 45233
$ws.Range('E49').Value = 10
$ws.Range('F49').Value = 300000000
$ws.Range('G49').Value = 'Espárragos'
$ws.Range('H49').Value = 'Sin especificar'
$ws.Range('I49').Value = 'Primera'
$ws.Range('J49').Value = 500
$ws.Range('K49').Value = 2000
$ws.Range('L49').Value = 2000
$ws.Range('M49').Value = 2000
$ws.Range('N49').Value = '$/kilo'
$ws.Range('O49').Value = 'Provincia de Linares'
$ws.Range('P49').Value = 2000
$ws.Range('Q49').Value = 1
$ws.Range('R49').Value = 'Hortaliza'
